$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '62.704.29'
Set-TextValue 'E2' '  -0.63%  '
Set-TextValue 'D3' '2.577.43'
Set-TextValue 'E3' '  +1.02%  '
Set-TextValue 'E4' '  +0.10%  '
Set-TextValue 'D5' '582.31'
Set-TextValue 'E5' '  -0.62%  '
Set-TextValue 'D6' '145.10'
Set-TextValue 'E6' '  -1.49%  '
Set-TextValue 'E7' '  +0.10%  '
Set-TextValue 'D8' '0.592'
Set-TextValue 'E8' '  +1.27%  '
Set-TextValue 'D9' '0.107'
Set-TextValue 'E9' '  +0.92%  '
Set-TextValue 'D10' '5.61'
Set-TextValue 'E10' '  +0.94%  '
Set-TextValue 'D12' '0.352'
Set-TextValue 'E12' '  -0.91%  '
Set-TextValue 'D13' '27.08'
Set-TextValue 'E13' '  -1.67%  '
Set-TextValue 'D14' '3.046.62'
Set-TextValue 'E14' '  +1.39%  '
Set-TextValue 'D15' '62.668.34'
Set-TextValue 'E15' '  -0.54%  '
Set-TextValue 'D16' '0.0000145'
Set-TextValue 'E16' '  +1.16%  '
Set-TextValue 'D17' '2.582.46'
Set-TextValue 'E17' '  +1.17%  '
Set-TextValue 'D18' '11.23'
Set-TextValue 'E18' '  -1.29%  '
Set-TextValue 'D19' '339.92'
Set-TextValue 'E19' '  +0.83%  '
Set-TextValue 'D20' '4.37'
Set-TextValue 'E20' '  +0.68%  '
Set-TextValue 'D21' '6.65'
Set-TextValue 'E21' '  -1.85%  '
Set-TextValue 'E22' '  -0.05%  '
Set-TextValue 'D23' '67.22'
Set-TextValue 'E23' '  +2.16%  '
Set-TextValue 'D24' '2.706.53'
Set-TextValue 'E24' '  +1.09%  '
Set-TextValue 'E25' '  -2.25%  '
Set-TextValue 'D26' '1.58'
Set-TextValue 'E26' '  -2.70%  '
Set-TextValue 'D27' '0.999'
Set-TextValue 'E27' '  -0.25%  '
Set-TextValue 'E28' '  +1.05%  '
Set-TextValue 'E29' '  -2.27%  '
Set-TextValue 'D30' '8.26'
Set-TextValue 'E30' '  -1.62%  '
Set-TextValue 'D31' '1.91'
Set-TextValue 'E31' '  -3.74%  '
Set-TextValue 'B32' 'PEPE'
Set-TextValue 'C32' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D32' '0.0₃0814'
Set-TextValue 'E32' '  -0.49%  '
Set-TextValue 'B33' 'Bittensor'
Set-TextValue 'C33' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D33' '458.46'
Set-TextValue 'E33' '  +9.00%  '
Set-TextValue 'D34' '176.43'
Set-TextValue 'E34' '  -1.22%  '
Set-TextValue 'E35' '  +3.10%  '
Set-TextValue 'E36' '  +0.06%  '
Set-TextValue 'D37' '0.397'
Set-TextValue 'E37' '  -1.14%  '
Set-TextValue 'D38' '18.92'
Set-TextValue 'E38' '  -1.14%  '
Set-TextValue 'D39' '4.50'
Set-TextValue 'E39' '  +2.97%  '
Set-TextValue 'D41' '1.69'
Set-TextValue 'E41' '  -3.40%  '
Set-TextValue 'D42' '157.94'
Set-TextValue 'E42' '  +5.02%  '
Set-TextValue 'D43' '3.73'
Set-TextValue 'E43' '  -1.72%  '
Set-TextValue 'B44' 'InjectiveProtocol'
Set-TextValue 'C44' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D44' '21.01'
Set-TextValue 'E44' '  +0.57%  '
Set-TextValue 'B45' 'Mantle'
Set-TextValue 'C45' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D45' '0.628'
Set-TextValue 'E45' '  +4.00%  '
Set-TextValue 'D46' '0.0537'
Set-TextValue 'E46' '  -0.95%  '
Set-TextValue 'D47' '0.0964'
Set-TextValue 'E47' '  -1.00%  '
Set-TextValue 'D48' '0.0235'
Set-TextValue 'E48' '  -1.21%  '
Set-TextValue 'D49' '18.18'
Set-TextValue 'E49' '  -0.83%  '
Set-TextValue 'B50' 'WhiteBITCoin'
Set-TextValue 'C50' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D50' '11.42'
Set-TextValue 'E50' '  +0.97%  '
Set-TextValue 'B51' 'dogwifhat'
Set-TextValue 'C51' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D51' '1.70'
Set-TextValue 'E51' '  -0.97%  '
